# PLAN 1.4 - increase merchant diversity in DNB Mastercard demo data
# Rewrites the April-2025 transaction list with a new, more diverse set of
# merchants, renames the worksheet, widens the datetime number format and
# extends the sheet with 5 additional transaction rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet (was "DNB Mastercard Demo") -----------------
$ws.Name = "Sheet"

# --- Row 2: date changes, merchant/amount stay the same ----------------
$ws.Range("A2").Value = 45777

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = 45776
$ws.Range("B3").Value = "KIWI STORO"
$ws.Range("F3").Value = 412.3

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = 45775
$ws.Range("B4").Value = "JOKER FROGNER"
$ws.Range("F4").Value = 234

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = 45773
$ws.Range("B5").Value = "NETFLIX.COM"
$ws.Range("F5").Value = 179

# --- Row 6 ---------------------------------------------------------------
$ws.Range("A6").Value = 45771
$ws.Range("B6").Value = "ODA.COM"
$ws.Range("F6").Value = 1567

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = 45769
$ws.Range("B7").Value = "VINMONOPOLET AKER BRYGGE"
$ws.Range("F7").Value = 456

# --- Row 8 ---------------------------------------------------------------
$ws.Range("A8").Value = 45767
$ws.Range("B8").Value = "FOODORA"
$ws.Range("F8").Value = 289

# --- Row 9 ---------------------------------------------------------------
$ws.Range("A9").Value = 45765
$ws.Range("B9").Value = "STARBUCKS OSLO S"
$ws.Range("F9").Value = 95

# --- Row 10 --------------------------------------------------------------
$ws.Range("A10").Value = 45764
$ws.Range("B10").Value = "WOLT"
$ws.Range("F10").Value = 345

# --- Row 11 --------------------------------------------------------------
$ws.Range("A11").Value = 45762
$ws.Range("B11").Value = "GITHUB.COM"
$ws.Range("F11").Value = 129

# --- Row 12 --------------------------------------------------------------
$ws.Range("A12").Value = 45761
$ws.Range("B12").Value = "JUST EAT"
$ws.Range("F12").Value = 267

# --- Row 13: was "Innbetaling" (E13=15000) -> now a regular purchase ----
$ws.Range("A13").Value = 45759
$ws.Range("B13").Value = "MENY CC VEST"
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = 534.6

# --- Row 14 --------------------------------------------------------------
$ws.Range("A14").Value = 45757
$ws.Range("B14").Value = "POWER LAMBERTSETER"
$ws.Range("F14").Value = 1899

# --- New rows 15-19 --------------------------------------------------------
$ws.Range("A15").Value = 45755
$ws.Range("B15").Value = "JACK & JONES OSLO CITY"
$ws.Range("F15").Value = 1199

$ws.Range("A16").Value = 45753
$ws.Range("B16").Value = "G-SPORT STORO"
$ws.Range("F16").Value = 1599

$ws.Range("A17").Value = 45752
$ws.Range("B17").Value = "Innbetaling"
$ws.Range("E17").Value = 15000

$ws.Range("A18").Value = 45750
$ws.Range("B18").Value = "COOP PRIX SAGENE"
$ws.Range("F18").Value = 534.2

$ws.Range("A19").Value = 45748
$ws.Range("B19").Value = "LINDEX OSLO CITY"
$ws.Range("F19").Value = 699

# Apply the date number format to the newly populated date cells too
$ws.Range("A2:A19").NumberFormat = "yyyy-mm-dd h:mm:ss"
